$d = $word.ActiveDocument

# 1. "Climatology applied" -> "Climatology" (Heading3)
$d.Content.Find.Execute("Climatology applied", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Climatology", 2)

# 2. "Ativação: 01/01/2024" -> "Ativação: 01/01/2025"
$d.Content.Find.Execute("Ativação: 01/01/2024", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ativação: 01/01/2025", 2)

# 3. Portuguese "Programa" paragraph: reorder the sentences around "Mudanças Climáticas"
$oldPt = "Caracterização Climática da Terra e do território brasileiro. Clima e suas relações com saúde, recursos hídricos, energia, agricultura. Previsão climática e modelos climáticos. Evolução do clima da Terra e Mudanças Climáticas.A disciplina pode contar"
$newPt = "Caracterização Climática da Terra e do território brasileiro. Evolução do clima da Terra e Mudanças Climáticas. Mudanças Climáticas e suas relações com saúde, recursos hídricos, energia, agricultura. A disciplina pode contar"
$d.Content.Find.Execute($oldPt, $true, $false, $false, $false, $false,
                         $true, 1, $false, $newPt, 2)

# 4. English "Programa" paragraph: add an extra space after "oceans's".
#    NOTE: this text contains a straight apostrophe ('); Find.Execute's
#    ReplaceWith argument goes through AutoFormat "smart quotes" and would
#    turn it into a curly apostrophe (’). Use a plain Range.Text assignment
#    instead, which performs a literal, un-autocorrected text replacement.
$rng = $d.Content
$rng.Find.Execute("oceans's general circulation")
if ($rng.Find.Found) {
    $rng.Text = "oceans's  general circulation"
}
